# Generate Report for Handback
# Refresh the timestamps recorded on the handback status report:
#  - Overview!G2            "Latest HO Xliff Generate Date" for the e57d00a9 file
#  - zh-cn!H2 / zh-cn!K2     "Correspond Handoff/Handback Datetime" for the e57d00a9 file (zh-cn)
#  - de-de!K2                "Correspond Handback Datetime" for the e57d00a9 file (de-de)

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-19 17:09:20"

$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-19 17:09:16"
$zhcn.Range("K2").Value = "2016-08-19 17:09:33"

$dede = $wb.Sheets.Item("de-de")
$dede.Range("K2").Value = "2016-08-19 17:09:40"
